# "better program creation to excel"
#
# 1. Set explicit custom column widths on all three sheets (program,
#    structures, sections).
# 2. Strip the placeholder empty inlineStr cells that used to pad out every
#    column of row 2+ (only cells that actually carry data survive).
# 3. On "sections", re-derive/re-order the per-layer country rows so they
#    line up layer-by-layer (US/Canada pairs per INSPER_ID_PRE) instead of
#    all-US-then-all-Canada.

# Excel's ColumnWidth (character units) differs from the raw OOXML <col
# width="..."> by a constant padding offset (~0.8333 for the default
# Calibri 11 font metrics used by this workbook). Subtract it so the saved
# width attribute lands on the exact integer values we want.
$colWidthPad = 0.8333333333333333

function Set-ColWidths {
    param($ws, [double[]]$widths)
    for ($i = 0; $i -lt $widths.Length; $i++) {
        $ws.Columns.Item($i + 1).ColumnWidth = $widths[$i] - $colWidthPad
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: program
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("program")

Set-ColWidths $ws1 @(15,28,12,14,19,16,25,27,29,31,25,27,25,36)

$ws1.Range("C2:D2").ClearContents()
$ws1.Range("F2:H2").ClearContents()
$ws1.Range("K2:N2").ClearContents()

# ---------------------------------------------------------------------
# Sheet 2: structures
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("structures")

Set-ColWidths $ws2 @(15,17,26,27,16,23,20,15,16,17,25,16,23,26,30,31,34,23,22,28,16)

$ws2.Range("B2:B4").ClearContents()
$ws2.Range("D2:D4").ClearContents()
$ws2.Range("F2:G4").ClearContents()
$ws2.Range("J2:Q4").ClearContents()
$ws2.Range("S2:U4").ClearContents()

# ---------------------------------------------------------------------
# Sheet 3: sections
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("sections")

Set-ColWidths $ws3 @(14,15,12,17,15,18,23,19,18,15,14,27,27,27,25,10,11,19,22,10,22,15,13,15,10,27,28,18,21,17,28,15,16,18,23,17,13,22)

$ws3.Range("C2:D7").ClearContents()
$ws3.Range("F2:H7").ClearContents()
$ws3.Range("J2:P7").ClearContents()
$ws3.Range("R2:R7").ClearContents()
$ws3.Range("T2:AG7").ClearContents()
$ws3.Range("AI2:AL7").ClearContents()

# Re-order the US/Canada rows so they group by layer (INSPER_ID_PRE):
# layer 1 (US, Canada), layer 2 (US, Canada), layer 3 (US, Canada) -- and
# refresh the Q/S/AH figures to match each new row's layer.
$rows = @(
    @{ E = 1; I = "United States"; Q = 8750000;  S = 3000000;  AH = 0.1    },
    @{ E = 1; I = "Canada";        Q = 8750000;  S = 3000000;  AH = 0.1    },
    @{ E = 2; I = "United States"; Q = 10000000; S = 11750000; AH = 0.1    },
    @{ E = 2; I = "Canada";        Q = 10000000; S = 11750000; AH = 0.1    },
    @{ E = 3; I = "United States"; Q = 23250000; S = 21750000; AH = 0.0979 },
    @{ E = 3; I = "Canada";        Q = 23250000; S = 21750000; AH = 0.0979 }
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    $ws3.Range("E$r").Value = $data.E
    $ws3.Range("I$r").Value = $data.I
    $ws3.Range("Q$r").Value = $data.Q
    $ws3.Range("S$r").Value = $data.S
    $ws3.Range("AH$r").Value = $data.AH
}
